$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 965; rows 965..1013 shift down to 966..1014.
$ws.Rows.Item(965).Insert()

# Populate the new row 965. Columns A,B,C,E,F,G,H,I,J,K,R keep the same
# values the (now shifted) row below it had before the edit; columns
# D,L,M,N,O,P,Q,S,T carry the new reading.
$ws.Cells.Item(965, 1).Value = 7
$ws.Cells.Item(965, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(965, 3).Value = "Ñuble"
$ws.Cells.Item(965, 4).Value = 45041
$ws.Cells.Item(965, 5).Value = 16
$ws.Cells.Item(965, 6).Value = "Fruta"
$ws.Cells.Item(965, 7).Value = 100102
$ws.Cells.Item(965, 8).Value = "Cítricos"
$ws.Cells.Item(965, 9).Value = 100102003
$ws.Cells.Item(965, 10).Value = "Limón"
$ws.Cells.Item(965, 11).Value = "Sin especificar"
$ws.Cells.Item(965, 12).Value = "1a plateado"
$ws.Cells.Item(965, 13).Value = 270
$ws.Cells.Item(965, 14).Value = 23000
$ws.Cells.Item(965, 15).Value = 24000
$ws.Cells.Item(965, 16).Value = 23556
$ws.Cells.Item(965, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(965, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(965, 19).Value = 1309
$ws.Cells.Item(965, 20).Value = 18
